$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Export" sheet is a Conta/Nome/Saldo listing sorted by Saldo (descending).
#
# 1) Add a new record for account 005995120 / Erik with a Saldo of 16250.
#    That balance sorts immediately above account 004218542 (Jose, 12760),
#    so insert a new row right before it and fill in the three columns.
$anchor = $ws.Cells.Find("004218542")
$insertRow = $anchor.Row

$ws.Rows($insertRow).Insert()

# Write the account number as text so the leading zeros are preserved, then
# strip the formatting Excel applies for a text-looking-like-a-number entry
# so the cell matches the plain (unstyled) look of the rest of the column.
$ws.Cells.Item($insertRow, 1).Value = "'005995120"
$ws.Cells.Item($insertRow, 1).ClearFormats()
$ws.Cells.Item($insertRow, 2).Value = "Erik"
$ws.Cells.Item($insertRow, 3).Value = 16250

# 2) Remove the existing record for account 005171652 / Bruno (Saldo 106.73).
$removed = $ws.Cells.Find("005171652")
$ws.Rows($removed.Row).Delete()
